$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold font, border, centered/top alignment) from H1 to I1 and J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set header values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Set I0/IF data values for rows 2-73
$data = @(
    @(2, 8, 8),
    @(3, 7, 7),
    @(4, 7, 7),
    @(5, 9, 9),
    @(6, 8, 8),
    @(7, 8, 8),
    @(8, 8, 9),
    @(9, 9, 9),
    @(10, 9, 9),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 9, 9),
    @(14, 9, 9),
    @(15, 9, 9),
    @(16, 7, 7),
    @(17, 9, 9),
    @(18, 9, 9),
    @(19, 9, 9),
    @(20, 9, 9),
    @(21, 9, 9),
    @(22, 9, 9),
    @(23, 9, 9),
    @(24, 9, 9),
    @(25, 9, 9),
    @(26, 10, 10),
    @(27, 9, 9),
    @(28, 9, 9),
    @(29, 9, 10),
    @(30, 9, 9),
    @(31, 9, 9),
    @(32, 9, 9),
    @(33, 9, 9),
    @(34, 9, 9),
    @(35, 9, 9),
    @(36, 9, 9),
    @(37, 8, 8),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 8, 9),
    @(41, 8, 9),
    @(42, 9, 9),
    @(43, 8, 9),
    @(44, 9, 9),
    @(45, 8, 8),
    @(46, 8, 9),
    @(47, 8, 10),
    @(48, 8, 9),
    @(49, 9, 9),
    @(50, 9, 9),
    @(51, 9, 9),
    @(52, 8, 9),
    @(53, 9, 9),
    @(54, 8, 8),
    @(55, 9, 9),
    @(56, 8, 9),
    @(57, 8, 9),
    @(58, 9, 9),
    @(59, 9, 10),
    @(60, 9, 9),
    @(61, 9, 9),
    @(62, 9, 9),
    @(63, 8, 9),
    @(64, 8, 9),
    @(65, 8, 8),
    @(66, 9, 9),
    @(67, 9, 9),
    @(68, 8, 8),
    @(69, 6, 6),
    @(70, 8, 8),
    @(71, 3, 4),
    @(72, 4, 4),
    @(73, 4, 5),
)

foreach ($entry in $data) {
    $rowNum = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($rowNum, 9).Value = $iVal   # column I
    $ws.Cells.Item($rowNum, 10).Value = $jVal  # column J
}

$wb.Application.CutCopyMode = $false
